$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.126.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.18%  "
$ws.Range("D3").Value = "'3.083.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.96%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'554.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'144.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.79%  "
$ws.Range("D8").Value = "'3.078.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.92%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'6.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +18.53%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.30%  "
$ws.Range("D12").Value = "'0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.80%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").Value = "'35.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.25%  "
$ws.Range("D15").Value = "'3.516.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'63.669.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "'3.056.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "'6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").Value = "'481.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("D21").Value = "'13.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.59%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.92%  "
$ws.Range("D23").Value = "'7.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.29%  "
$ws.Range("D24").Value = "'13.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.39%  "
$ws.Range("D25").Value = "'80.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'2.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").Value = "'7.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.07%  "
$ws.Range("D29").Value = "'2.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.54%  "
$ws.Range("D30").Value = "'0.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "'26.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "'2.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.62%  "
$ws.Range("D34").Value = "'5.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.92%  "
$ws.Range("D35").Value = "'55.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").Value = "'6.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.09%  "
$ws.Range("D37").Value = "'461.57"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.40%  "
$ws.Range("D38").Value = "'0.0833"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.83%  "
$ws.Range("D39").Value = "'0.0404"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.37%  "
$ws.Range("D40").Value = "'3.000.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("D41").Value = "'0.118"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").Value = "'8.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("D43").Value = "'2.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.74%  "
$ws.Range("D44").Value = "'27.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.45%  "
$ws.Range("D45").Value = "'0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.16%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'2.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.52%  "
$ws.Range("E48").Value = "  +4.51%  "
$ws.Range("D49").Value = "'0.0₃0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.58%  "
$ws.Range("D50").Value = "'116.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'2.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.47%  "
